# Arbeitsaufteilung.xlsx - "Added responsibility for working packages"
#
# - Fills in the work-package names (column A) for the existing task rows
# - Highlights (yellow) the "e." estimate-column cell of whoever is
#   responsible for each work package
# - Inserts two extra blank rows before the total/sum rows
# - Widens columns A and B so the longer work-package text is readable
# - Adjusts a couple of row heights (wrapped / manually resized rows)
# - Extends the print area to cover the two new rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column widths: A and B both become wide (33) ----------------------
$ws.Columns.Item(1).ColumnWidth = 32.1
$ws.Columns.Item(2).ColumnWidth = 32.1

# --- Work package names (column A, rows 3-20) ---------------------------
$ws.Range("A3").Value  = "UML Klassendiagramm"
$ws.Range("A4").Value  = "UML Aktivitätsdiagramm"
$ws.Range("A5").Value  = "UML Use-Case"
$ws.Range("A6").Value  = "UML überprüfen"
$ws.Range("A7").Value  = "Analytics-Server implementieren"
$ws.Range("A8").Value  = "Billing-Server implementieren"
$ws.Range("A9").Value  = "Management-Client impl."
$ws.Range("A10").Value = "Testing Component impl."
$ws.Range("A11").Value = "Model-Klassen (Events, Bill, Steps)"
$ws.Range("A12").Value = "File-Persistence"
$ws.Range("A13").Value = "Refactoring old Source"
$ws.Range("A14").Value = "RMI-Verbindungen implementieren"
$ws.Range("A15").Value = "RMI-Verbindungen testen"
$ws.Range("A16").Value = "Analytics Unit testen"
$ws.Range("A17").Value = "Billing Unit testen"
$ws.Range("A18").Value = "Management-Client Unit testen"
$ws.Range("A19").Value = "Testing Component Unit testen"
$ws.Range("A20").Value = "Protokoll"

# --- Row height tweaks ---------------------------------------------------
# Row 7 text wraps onto two lines at the new column width.
$ws.Rows.Item(7).RowHeight = 30.75
# Row 14 got a small manual height bump.
$ws.Rows.Item(14).RowHeight = 17.25

# --- Responsibility highlighting (yellow fill on the "e." cell) --------
# Columns: B/C=Reichmann, D/E=Krepela, F/G=Lipovits, H/I=Tattyrek, J/K=Traxler
$yellow = 65535  # RGB(255,255,0)

$responsible = @{
  "D3"  = $true; "J3"  = $true
  "F4"  = $true
  "B5"  = $true
  "H6"  = $true
  "H7"  = $true
  "D8"  = $true; "J8"  = $true
  "F9"  = $true
  "F10" = $true
  "H11" = $true
  "B13" = $true
  "B14" = $true; "J14" = $true
  "J15" = $true
  "B16" = $true; "H16" = $true
  "D17" = $true
  "F18" = $true
  "F19" = $true
  "B20" = $true; "D20" = $true; "F20" = $true; "H20" = $true; "J20" = $true
}

foreach ($addr in $responsible.Keys) {
  $ws.Range($addr).Interior.Color = $yellow
}

# --- Insert two blank rows before the total row -------------------------
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()

# Give the two new blank rows the same look as the row above them.
$ws.Range("A21:K21").Copy()
$ws.Range("A22:K22").PasteSpecial(-4122)
$ws.Range("A23:K23").PasteSpecial(-4122)

# --- Print area now covers the two extra rows ---------------------------
$ws.PageSetup.PrintArea = "`$A`$1:`$I`$24"

# --- Re-create the split window / selection state -----------------------
$ws.Activate()
$ws.Range("A15").Select()
$excel.ActiveWindow.SplitColumn = 3
$excel.ActiveWindow.SplitRow = 0
$ws.Range("J8").Select()

"Done"
